$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Receipts" table (header row 7 / data row 8) gets a new "Date" column
# inserted right after "ReceiptUniqueId" (column B), pushing Title/Total/
# PaidBy/FK - AccountID one column to the right (C:F -> D:G).
#
# Shift the existing cells rightward first (right-to-left so nothing is
# clobbered before it's copied), using .Copy() so each cell keeps its
# original formatting/style attached instead of just the raw value.

# Row 7 (header): C..F -> D..G
$ws.Range("F7").Copy($ws.Range("G7"))
$ws.Range("E7").Copy($ws.Range("F7"))
$ws.Range("D7").Copy($ws.Range("E7"))
$ws.Range("C7").Copy($ws.Range("D7"))
$ws.Range("C7").Style = "Normal"
$ws.Range("C7").Value = "Date"

# Row 8 (data): C..F -> D..G
$ws.Range("F8").Copy($ws.Range("G8"))
$ws.Range("E8").Copy($ws.Range("F8"))
$ws.Range("D8").Copy($ws.Range("E8"))
$ws.Range("C8").Copy($ws.Range("D8"))
$ws.Range("C8").Style = "Normal"

# New Date cell: 31 May 2020 as an Excel serial date, formatted with the
# built-in short-date number format (numFmtId 14).
$ws.Range("C8").Value = 43982
$ws.Range("C8").NumberFormat = "mm-dd-yy"

# Move the active selection to B17 (was E16).
$ws.Range("B17").Select()
